# Add the 2024/11/03 data column (BD) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell BD1: date label stored as literal text, same style as BC1 ---
$ws.Cells.Item(1, 56).NumberFormat = "@"
$ws.Cells.Item(1, 56).Value = "2024/11/03"
$ws.Cells.Item(1, 55).Copy()
$ws.Cells.Item(1, 56).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New column width, matching the other data columns (stored width 12) ---
$ws.Columns.Item(56).ColumnWidth = 11.166666666666666

# --- Data values BD2:BD53 ---
$ws.Cells.Item(2, 56).Value = 149.8
$ws.Cells.Item(3, 56).Value = 142.9
$ws.Cells.Item(4, 56).Value = 218
$ws.Cells.Item(5, 56).Value = 217.4
$ws.Cells.Item(6, 56).Value = 155.5
$ws.Cells.Item(7, 56).Value = 179.4
$ws.Cells.Item(8, 56).Value = 129.1
$ws.Cells.Item(9, 56).Value = 152.2
$ws.Cells.Item(10, 56).Value = 135.3
$ws.Cells.Item(11, 56).Value = 147.4
$ws.Cells.Item(12, 56).Value = 142.6
$ws.Cells.Item(13, 56).Value = 166.5
$ws.Cells.Item(14, 56).Value = 179.5
$ws.Cells.Item(15, 56).Value = 154.8
$ws.Cells.Item(16, 56).Value = 150.8
$ws.Cells.Item(17, 56).Value = 105.6
$ws.Cells.Item(18, 56).Value = 123.3
$ws.Cells.Item(19, 56).Value = 183
$ws.Cells.Item(20, 56).Value = 155.4
$ws.Cells.Item(21, 56).Value = 141.6
$ws.Cells.Item(22, 56).Value = 127.3
$ws.Cells.Item(23, 56).Value = 156.3
$ws.Cells.Item(24, 56).Value = 150.7
$ws.Cells.Item(25, 56).Value = 132.7
$ws.Cells.Item(26, 56).Value = 165.2
$ws.Cells.Item(27, 56).Value = 188.7
$ws.Cells.Item(28, 56).Value = 195.5
$ws.Cells.Item(29, 56).Value = 258.9
$ws.Cells.Item(30, 56).Value = 140.8
$ws.Cells.Item(31, 56).Value = 194.8
$ws.Cells.Item(32, 56).Value = 144.9
$ws.Cells.Item(33, 56).Value = 216.3
$ws.Cells.Item(34, 56).Value = 140.1
$ws.Cells.Item(35, 56).Value = 173.6
$ws.Cells.Item(36, 56).Value = 159.4
$ws.Cells.Item(37, 56).Value = 131.3
$ws.Cells.Item(38, 56).Value = 122.7
$ws.Cells.Item(39, 56).Value = 183.2
$ws.Cells.Item(40, 56).Value = 152.2
$ws.Cells.Item(41, 56).Value = 189.1
$ws.Cells.Item(42, 56).Value = 119.3
$ws.Cells.Item(43, 56).Value = 137.7
$ws.Cells.Item(44, 56).Value = 153.9
$ws.Cells.Item(45, 56).Value = 148.2
$ws.Cells.Item(46, 56).Value = 139.3
$ws.Cells.Item(47, 56).Value = 181.4
$ws.Cells.Item(48, 56).Value = 162.3
$ws.Cells.Item(49, 56).Value = 149.1
$ws.Cells.Item(50, 56).Value = 170.2
$ws.Cells.Item(51, 56).Value = 152.7
$ws.Cells.Item(52, 56).Value = 167.3
$ws.Cells.Item(53, 56).Value = 146.9

# --- Default highlight style (no fill) for the whole data range ---
$ws.Cells.Item(3, 55).Copy()
$ws.Range("BD2:BD53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Override style for yellow (value < 125) ---
$ws.Cells.Item(4, 55).Copy()
$ws.Cells.Item(17, 56).PasteSpecial(-4122)
$ws.Cells.Item(18, 56).PasteSpecial(-4122)
$ws.Cells.Item(38, 56).PasteSpecial(-4122)
$ws.Cells.Item(42, 56).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Override style for blue (125 <= value < 140) ---
$ws.Cells.Item(2, 55).Copy()
$ws.Cells.Item(8, 56).PasteSpecial(-4122)
$ws.Cells.Item(10, 56).PasteSpecial(-4122)
$ws.Cells.Item(22, 56).PasteSpecial(-4122)
$ws.Cells.Item(25, 56).PasteSpecial(-4122)
$ws.Cells.Item(37, 56).PasteSpecial(-4122)
$ws.Cells.Item(43, 56).PasteSpecial(-4122)
$ws.Cells.Item(46, 56).PasteSpecial(-4122)
$excel.CutCopyMode = $false

